# Generate Report for Handoff
#
# The file "96e78804-1ed0-4da2-b7d2-3f6c213bc972.md" has moved from
# "Handed back: in sync with en-US" to "Ready for handoff" for both the
# zh-cn and de-de locales, and its handoff timestamps were refreshed.
#
# Row 3 on every worksheet corresponds to this file.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"   # zh-cn status
$overview.Range("C3").Value = "Ready for handoff"   # de-de status
$overview.Range("D3").Value = "2016-51-13 08:51:12" # Latest Handoff Date

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"        # Status
$zhcn.Range("E3").Value = "2016-03-13 08:51:09"      # Latest Handoff Datetime

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"        # Status
$dede.Range("E3").Value = "2016-03-13 08:51:12"      # Latest Handoff Datetime
